$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.705.12"
$ws.Range("E2").Value = "  -3.04%  "
$ws.Range("D3").Value = "2.098.09"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "'343.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.93%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "'0.5137"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("D8").Value = "'0.4403"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'53.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "'0.09173"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "'24.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "2.096.56"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "'6.752"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "'8.165"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").Value = "'1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'20.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.41%  "
$ws.Range("D20").Value = "'0.06641"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "'1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "'6.175"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").Value = "29.745.62"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").Value = "'12.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "'2.292"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.97%  "
$ws.Range("D26").Value = "2.345.03"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Value = "'21.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("D28").Value = "'162.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").Value = "'2.507"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").Value = "'132.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("D31").Value = "'1.131"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.49%  "
$ws.Range("D32").Value = "'0.1047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("D33").Value = "'1.636"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "'6.148"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("D35").Value = "'3.971"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").Value = "'6.058"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "'10.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").Value = "'0.02561"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").Value = "'0.06693"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").Value = "'12.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").Value = "'0.2210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("D43").Value = "'1.300"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").Value = "'0.6659"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").Value = "'14.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("D46").Value = "'2.297"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Value = "'3.607"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.51%  "
$ws.Range("D48").Value = "'1.218"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'81.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000334"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.00%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "'1.161"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.83%  "
